$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only "input_Name" (currently in column I) should survive, becoming column A.
# Move column I to the front (column A) by cutting it and inserting it at position 1.
$ws.Columns.Item(9).Cut()
$ws.Columns.Item(1).Insert()

# Remove all other (now shifted) columns B:M, leaving only the moved column A.
$ws.Range("B1:M2").EntireColumn.Delete()

Write-Host "done"
